# Update available winget packages: mark column E ("gh/microsoft/winget-pkgs")
# with "x" for packages that are now available via winget.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$rows = @(4, 9, 12, 18, 31, 38, 41, 46, 50, 52, 55, 62, 67, 71, 74, 87, 92, 94, 95, 96, 101, 107, 108, 109, 113)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "x"
}

# Restore the selection/frozen-pane view state to match the author's session.
$ws.Activate()
($excel.ActiveWindow.FreezePanes = $true) | Out-Null
$ws.Range("E116").Select() | Out-Null
